$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2049689440993789
$ws.Range("C2").Value = 0.5186335403726708
$ws.Range("J2").Value = 0.01863354037267081
$ws.Range("P2").Value = 0.1521739130434783
$ws.Range("S2").Value = 0.1055900621118012
$ws.Range("B3").Value = 0.005988023952095809
$ws.Range("J3").Value = 0.02395209580838323
$ws.Range("P3").Value = 0.7305389221556886
$ws.Range("S3").Value = 0.2395209580838323
$ws.Range("J4").Value = 0.08
$ws.Range("P4").Value = 0.66
$ws.Range("S4").Value = 0.26
$ws.Range("B6").Value = 0.07281553398058252
$ws.Range("D6").Value = 0.02427184466019417
$ws.Range("F6").Value = 0.03883495145631068
$ws.Range("J6").Value = 0.3155339805825243
$ws.Range("O6").Value = 0.009708737864077669
$ws.Range("Q6").Value = 0.1699029126213592
$ws.Range("R6").Value = 0.03883495145631068
$ws.Range("S6").Value = 0.3300970873786408
$ws.Range("B7").Value = 0.1242937853107345
$ws.Range("D7").Value = 0.03389830508474576
$ws.Range("F7").Value = 0.01694915254237288
$ws.Range("J7").Value = 0.1299435028248588
$ws.Range("O7").Value = 0.01694915254237288
$ws.Range("Q7").Value = 0.1468926553672316
$ws.Range("R7").Value = 0.06214689265536723
$ws.Range("S7").Value = 0.4689265536723164
$ws.Range("B8").Value = 0.09070796460176991
$ws.Range("D8").Value = 0.02654867256637168
$ws.Range("E8").Value = 0.002212389380530973
$ws.Range("F8").Value = 0.05088495575221239
$ws.Range("J8").Value = 0.1106194690265487
$ws.Range("O8").Value = 0.01548672566371681
$ws.Range("Q8").Value = 0.1570796460176991
$ws.Range("R8").Value = 0.1150442477876106
$ws.Range("S8").Value = 0.4314159292035398
$ws.Range("B9").Value = 0.05825242718446602
$ws.Range("D9").Value = 0.02427184466019417
$ws.Range("F9").Value = 0.03398058252427184
$ws.Range("J9").Value = 0.116504854368932
$ws.Range("O9").Value = 0.01456310679611651
$ws.Range("Q9").Value = 0.2281553398058253
$ws.Range("R9").Value = 0.1262135922330097
$ws.Range("S9").Value = 0.3980582524271845
$ws.Range("B10").Value = 0.1169895678092399
$ws.Range("D10").Value = 0.01937406855439642
$ws.Range("F10").Value = 0.05812220566318927
$ws.Range("J10").Value = 0.139344262295082
$ws.Range("O10").Value = 0.01117734724292101
$ws.Range("Q10").Value = 0.1944858420268256
$ws.Range("R10").Value = 0.08196721311475409
$ws.Range("S10").Value = 0.3785394932935917
$ws.Range("G11").Value = 0.1331058020477816
$ws.Range("J11").Value = 0.1331058020477816
$ws.Range("K11").Value = 0.204778156996587
$ws.Range("L11").Value = 0.5085324232081911
$ws.Range("S11").Value = 0.0204778156996587
$ws.Range("G12").Value = 0.7161290322580646
$ws.Range("J12").Value = 0.2129032258064516
$ws.Range("K12").Value = 0.01290322580645161
$ws.Range("L12").Value = 0.03870967741935484
$ws.Range("S12").Value = 0.01935483870967742
$ws.Range("F13").Value = 0.02040816326530612
$ws.Range("G13").Value = 0.6326530612244898
$ws.Range("J13").Value = 0.2857142857142857
$ws.Range("S13").Value = 0.06122448979591837
$ws.Range("F15").Value = 0.03404255319148936
$ws.Range("H15").Value = 0.1659574468085106
$ws.Range("I15").Value = 0.06808510638297872
$ws.Range("J15").Value = 0.4
$ws.Range("K15").Value = 0.05106382978723404
$ws.Range("M15").Value = 0.008510638297872341
$ws.Range("O15").Value = 0.05531914893617021
$ws.Range("S15").Value = 0.2170212765957447
$ws.Range("F16").Value = 0.02030456852791878
$ws.Range("H16").Value = 0.1522842639593909
$ws.Range("I16").Value = 0.09137055837563451
$ws.Range("J16").Value = 0.3756345177664975
$ws.Range("K16").Value = 0.1116751269035533
$ws.Range("M16").Value = 0.02030456852791878
$ws.Range("O16").Value = 0.06091370558375635
$ws.Range("S16").Value = 0.16751269035533
$ws.Range("F17").Value = 0.006880733944954129
$ws.Range("H17").Value = 0.2018348623853211
$ws.Range("I17").Value = 0.08027522935779817
$ws.Range("J17").Value = 0.3761467889908257
$ws.Range("K17").Value = 0.1032110091743119
$ws.Range("M17").Value = 0.02293577981651376
$ws.Range("N17").Value = 0.002293577981651376
$ws.Range("O17").Value = 0.07798165137614679
$ws.Range("S17").Value = 0.1284403669724771
$ws.Range("F18").Value = 0.04411764705882353
$ws.Range("H18").Value = 0.1617647058823529
$ws.Range("I18").Value = 0.1127450980392157
$ws.Range("J18").Value = 0.3725490196078431
$ws.Range("K18").Value = 0.07843137254901961
$ws.Range("M18").Value = 0.01470588235294118
$ws.Range("N18").Value = 0.004901960784313725
$ws.Range("O18").Value = 0.07843137254901961
$ws.Range("S18").Value = 0.1323529411764706
$ws.Range("F19").Value = 0.02552552552552553
$ws.Range("H19").Value = 0.2004504504504505
$ws.Range("I19").Value = 0.08633633633633633
$ws.Range("J19").Value = 0.3798798798798799
$ws.Range("K19").Value = 0.09834834834834835
$ws.Range("M19").Value = 0.02177177177177177
$ws.Range("O19").Value = 0.07807807807807808
$ws.Range("S19").Value = 0.1096096096096096
